$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateRecipient")

$ws.Range("A4").Value = "TrialData"
$ws.Range("B4").Value = "'Recipient"
$ws.Range("C4").Value = "'91"

$ws.Range("A5").Value = "TrialData"
$ws.Range("B5").Value = "'Recipient"
$ws.Range("C5").Value = "'27"

$ws.Range("C6").Select()
